# Remove three duplicate/outlier weekly price rows from the Cilantro sheet.
# Original (pre-edit) row numbers 183, 175 and 137 are deleted. Deleting from
# the bottom up keeps the remaining row numbers valid as each deletion is
# applied, and Excel automatically shifts the rows below each deleted row
# upward and shrinks the worksheet's used range / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(183).Delete()
$ws.Rows(175).Delete()
$ws.Rows(137).Delete()
